# SMARTNODES.docx — apply translated/updated SmartNode services bullets,
# fix the bookmark id bump that came along with the edit.

$d = $word.ActiveDocument

# --- 1) "InstantPay (instant transactions)..." bullet -> new InstantPay copy, numId 1 -> 2
$p1 = $d.Paragraphs(5)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($xml1)

# --- 2) "Will have more services added later" -> SmartRewards bullet (numId stays 2)
$p2 = $d.Paragraphs(6)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)

# --- 3) bookmark "smarthosting" id bumped from 0 to 1 (paragraph text/format unchanged)
$bk = $d.Bookmarks("smarthosting")
$p3 = $bk.Range.Paragraphs(1)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="264" w:lineRule="atLeast"/><w:jc w:val="center"/><w:textAlignment w:val="baseline"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:caps/><w:spacing w:val="15"/><w:kern w:val="36"/><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="smarthosting"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica" w:eastAsia="Times New Roman"/><w:caps/><w:spacing w:val="15"/><w:kern w:val="36"/><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr><w:t xml:space="preserve">WHAT DO I NEED TO HOST A SMARTNODE?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xml3)

Write-Host "done"
